$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3346172571182251
$ws.Range("B1").Value = 0.8125601410865784
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.048319697380066
